$d = $word.ActiveDocument

# Author used for the tracked changes introduced below.
$word.UserName = "AlGord"

# --- 1. Move the "_GoBack" bookmark to the very start of the document ------
# (Word drops "_GoBack" at the last edit point; remove it from its old spot
#  next to "Our Background" and recreate it spanning the document's first
#  paragraph mark, same as a fresh edit session would.)
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
$startRange = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $startRange)

# --- 2. Track-change deletion in the "Customers" paragraph ------------------
# " ... have a need for a simple, one-stop-shop, document management
# solution." -> keep "... have a ne" + "ed for a simple, one-stop-shop",
# delete ", document management solution." (tracked).
$splitRange1 = $d.Content
$splitRange1.Find.Execute("ed for a simple, one-stop-shop")
$splitStart = $splitRange1.Start
$splitEnd = $splitRange1.End

# Force a clean run split at the two boundaries without generating a
# revision mark (toggle formatting on/off so Word cannot silently merge the
# surrounding runs back together).
$splitRangeA = $d.Range($splitStart, $splitEnd)
$splitRangeA.Font.Bold = $true
$splitRangeB = $d.Range($splitStart, $splitEnd)
$splitRangeB.Font.Bold = $false

$d.TrackRevisions = $true
$deleteRange = $d.Content
$deleteRange.Find.Execute(", document management solution.")
$deleteRange.Delete()

# --- 3. Track-change insertion in the "Our Background" paragraph ----------
# " ... Aspose's experience, stability and award winning technology." ->
# keep " ... Aspose's experience" and insert ", stability and award winning
# technology." as a tracked insertion.
$d.TrackRevisions = $true
$tailRange = $d.Content
$tailRange.Find.Execute(" experience, stability and award winning technology.")
$insertAt = $tailRange.Start + " experience".Length
$oldTailLen = $tailRange.End - $insertAt

$insPoint = $d.Range($insertAt, $insertAt)
$insPoint.InsertBefore(", stability and award winning technology.")

$newTailLen = ", stability and award winning technology.".Length
$oldTailStart = $insertAt + $newTailLen
$oldTailEnd = $oldTailStart + $oldTailLen

$d.TrackRevisions = $false
$oldTailRange = $d.Range($oldTailStart, $oldTailEnd)
$oldTailRange.Delete()
$d.TrackRevisions = $true
